$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.794.43'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '2.951.86'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'593.38"
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').Value = "'148.16"
$ws.Range('E6').Value = '  +2.07%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = "'0.509"
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '2.955.07'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = "'7.33"
$ws.Range('E10').Value = '  +5.26%  '
$ws.Range('D11').Value = "'0.152"
$ws.Range('E11').Value = '  +7.26%  '
$ws.Range('D12').Value = "'0.444"
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = "'0.0000240"
$ws.Range('E13').Value = '  +6.44%  '
$ws.Range('D14').Value = "'32.93"
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '3.444.26'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '62.789.34'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = "'6.73"
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '2.969.04'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('D20').Value = "'442.56"
$ws.Range('E20').Value = '  +2.27%  '
$ws.Range('D21').Value = "'13.50"
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = "'0.669"
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').Value = "'7.06"
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'81.56"
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D25').Value = "'11.16"
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('D26').Value = "'2.16"
$ws.Range('E26').Value = '  -2.15%  '
$ws.Range('D27').Value = "'11.70"
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = "'2.25"
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').Value = "'0.0000105"
$ws.Range('E30').Value = '  +19.51%  '
$ws.Range('D31').Value = "'7.26"
$ws.Range('E31').Value = '  +4.65%  '
$ws.Range('D32').Value = "'2.62"
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').Value = "'26.54"
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').Value = "'0.109"
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').Value = "'0.994"
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = "'3.21"
$ws.Range('E37').Value = '  +6.50%  '
$ws.Range('D38').Value = "'5.62"
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').Value = "'2.04"
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').Value = "'49.65"
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = "'8.53"
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').Value = "'0.118"
$ws.Range('E42').Value = '  -4.99%  '
$ws.Range('D43').Value = "'0.282"
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').Value = "'39.96"
$ws.Range('E44').Value = '  -6.56%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.701.39'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = "'135.49"
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('D47').Value = "'0.0340"
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = "'365.58"
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D50').Value = "'23.01"
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('E51').Value = '  -0.31%  '
